$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.904.83'
$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").Value = '2.444.04'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.31'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.58%  '
$ws.Range("E7").Value = '  +0.85%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.79'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.122'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.14'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.05'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '2.829.31'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '2.454.18'
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.838'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '45.766.42'
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.41'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  +2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.17'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.26%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.36'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.30%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '246.14'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.88'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.59'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.28'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.126'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.89'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.35'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0757'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.52'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("E39").Value = '  +0.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '127.33'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("E41").Value = '  +4.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.110'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.82'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0291'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("D45").Value = '1.954.02'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.95'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("E48").Value = '  +10.67%  '
$ws.Range("E49").Value = '  -5.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.96'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +6.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.25'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.78%  '
